# Sửa "Trần Văn Thành" -> "Văn-Thành" trong phân công chuyên môn phòng nhạc
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("phòng nhạc")

$ws.Range("E4").Value = "Âm nhạc-Văn-Thành-2B"
$ws.Range("F4").Value = "Âm nhạc-Văn-Thành-2C"
$ws.Range("G4").Value = "Âm nhạc-Văn-Thành-2Đ"
$ws.Range("G5").Value = "Âm nhạc-Văn-Thành-1C"
$ws.Range("E6").Value = "Âm nhạc-Văn-Thành-4A"
$ws.Range("F6").Value = "Âm nhạc-Văn-Thành-2D"
$ws.Range("G6").Value = "Âm nhạc-Văn-Thành-2A"
$ws.Range("E7").Value = "Âm nhạc-Văn-Thành-4B"
$ws.Range("F7").Value = "Âm nhạc-Văn-Thành-1D"
$ws.Range("G7").Value = "Âm nhạc-Văn-Thành-5B"
$ws.Range("E8").Value = "Âm nhạc-Văn-Thành-5A"
$ws.Range("C11").Value = "Âm nhạc-Văn-Thành-3A"
$ws.Range("D11").Value = "Âm nhạc-Văn-Thành-3C"
$ws.Range("D12").Value = "Âm nhạc-Văn-Thành-3B"

# Excel auto-recomputed "best fit" width for C:G after the shorter names were written
$ws.Columns("C:G").ColumnWidth = 23.8

Write-Host "done"
